$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows 175-189 with revised weekly values ----
$ws.Range("D175").Value = 44505
$ws.Range("M175").Value = 200
$ws.Range("N175").Value = 8000
$ws.Range("O175").Value = 8000
$ws.Range("P175").Value = 8000
$ws.Range("S175").Value = 1143
$ws.Range("D176").Value = 44505
$ws.Range("M176").Value = 200
$ws.Range("N176").Value = 6500
$ws.Range("O176").Value = 6500
$ws.Range("P176").Value = 6500
$ws.Range("S176").Value = 929
$ws.Range("D177").Value = 44505
$ws.Range("M177").Value = 100
$ws.Range("N177").Value = 5500
$ws.Range("O177").Value = 5500
$ws.Range("P177").Value = 5500
$ws.Range("S177").Value = 786
$ws.Range("D178").Value = 44168
$ws.Range("N178").Value = 10000
$ws.Range("O178").Value = 10000
$ws.Range("P178").Value = 10000
$ws.Range("S178").Value = 1429
$ws.Range("D179").Value = 44168
$ws.Range("N179").Value = 8000
$ws.Range("O179").Value = 8000
$ws.Range("P179").Value = 8000
$ws.Range("S179").Value = 1143
$ws.Range("D180").Value = 44168
$ws.Range("L180").Value = "Segunda"
$ws.Range("M180").Value = 50
$ws.Range("N180").Value = 7000
$ws.Range("O180").Value = 7000
$ws.Range("P180").Value = 7000
$ws.Range("S180").Value = 1000
$ws.Range("D181").Value = 44484
$ws.Range("L181").Value = "Especial"
$ws.Range("D182").Value = 44484
$ws.Range("L182").Value = "Primera"
$ws.Range("R182").Value = "Provincia de Melipilla"
$ws.Range("D183").Value = 44328
$ws.Range("L183").Value = "Especial"
$ws.Range("N183").Value = 15000
$ws.Range("O183").Value = 15000
$ws.Range("P183").Value = 15000
$ws.Range("R183").Value = "Provincia de Melipilla"
$ws.Range("S183").Value = 2143
$ws.Range("D184").Value = 44328
$ws.Range("L184").Value = "Primera"
$ws.Range("N184").Value = 12000
$ws.Range("O184").Value = 12000
$ws.Range("P184").Value = 12000
$ws.Range("R184").Value = "Provincia de Melipilla"
$ws.Range("S184").Value = 1714
$ws.Range("D185").Value = 44217
$ws.Range("L185").Value = "Especial"
$ws.Range("M185").Value = 100
$ws.Range("N185").Value = 9000
$ws.Range("O185").Value = 9000
$ws.Range("P185").Value = 9000
$ws.Range("R185").Value = "Región del Maule"
$ws.Range("S185").Value = 1286
$ws.Range("D186").Value = 44217
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 100
$ws.Range("N186").Value = 8000
$ws.Range("O186").Value = 8000
$ws.Range("P186").Value = 8000
$ws.Range("R186").Value = "Región del Maule"
$ws.Range("S186").Value = 1143
$ws.Range("D187").Value = 44217
$ws.Range("L187").Value = "Segunda"
$ws.Range("M187").Value = 100
$ws.Range("N187").Value = 7000
$ws.Range("O187").Value = 7000
$ws.Range("P187").Value = 7000
$ws.Range("R187").Value = "Región del Maule"
$ws.Range("S187").Value = 1000
$ws.Range("D188").Value = 44421
$ws.Range("M188").Value = 50
$ws.Range("N188").Value = 22000
$ws.Range("O188").Value = 22000
$ws.Range("P188").Value = 22000
$ws.Range("S188").Value = 3143
$ws.Range("D189").Value = 44433
$ws.Range("M189").Value = 50
$ws.Range("N189").Value = 24000
$ws.Range("O189").Value = 24000
$ws.Range("P189").Value = 24000
$ws.Range("S189").Value = 3429
# ---- Append 3 new data rows (190-192) for the latest week ----
# Row 190
$ws.Range("A190").Value = 11
$ws.Range("B190").Value = "Vega Monumental Concepción"
$ws.Range("C190").Value = "Bíobío"
$ws.Range("D190").Value = 44491
$ws.Range("E190").Value = 8
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100101
$ws.Range("H190").Value = "Berries"
$ws.Range("I190").Value = 100112025
$ws.Range("J190").Value = "Frutilla"
$ws.Range("K190").Value = "Sin especificar"
$ws.Range("L190").Value = "Especial"
$ws.Range("M190").Value = 270
$ws.Range("N190").Value = 9000
$ws.Range("O190").Value = 10000
$ws.Range("P190").Value = 9556
$ws.Range("Q190").Value = "$/bandeja 7 kilos"
$ws.Range("R190").Value = "Provincia de Melipilla"
$ws.Range("S190").Value = 1365
$ws.Range("T190").Value = 7
# Row 191
$ws.Range("A191").Value = 11
$ws.Range("B191").Value = "Vega Monumental Concepción"
$ws.Range("C191").Value = "Bíobío"
$ws.Range("D191").Value = 44491
$ws.Range("E191").Value = 8
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100101
$ws.Range("H191").Value = "Berries"
$ws.Range("I191").Value = 100112025
$ws.Range("J191").Value = "Frutilla"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 380
$ws.Range("N191").Value = 6500
$ws.Range("O191").Value = 7000
$ws.Range("P191").Value = 6737
$ws.Range("Q191").Value = "$/bandeja 7 kilos"
$ws.Range("R191").Value = "Provincia de Melipilla"
$ws.Range("S191").Value = 962
$ws.Range("T191").Value = 7
# Row 192
$ws.Range("A192").Value = 11
$ws.Range("B192").Value = "Vega Monumental Concepción"
$ws.Range("C192").Value = "Bíobío"
$ws.Range("D192").Value = 44491
$ws.Range("E192").Value = 8
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100101
$ws.Range("H192").Value = "Berries"
$ws.Range("I192").Value = 100112025
$ws.Range("J192").Value = "Frutilla"
$ws.Range("K192").Value = "Sin especificar"
$ws.Range("L192").Value = "Segunda"
$ws.Range("M192").Value = 200
$ws.Range("N192").Value = 5500
$ws.Range("O192").Value = 5500
$ws.Range("P192").Value = 5500
$ws.Range("Q192").Value = "$/bandeja 7 kilos"
$ws.Range("R192").Value = "Provincia de Melipilla"
$ws.Range("S192").Value = 786
$ws.Range("T192").Value = 7
$ws.Range("D190:D192").NumberFormat = "YYYY-MM-DD HH:MM:SS"
